# TI_Senior_Time_Tracking.xlsx - "Added time worked for Kevin"
#
# Fills in the previously-blank last row (row 76) of the time-tracking
# table with a new entry for Kevin Su: a work session on 2025-01-19 from
# 2:00 PM to 4:30 PM (2.5 hours) describing work on the mmWave Radar PCB
# schematic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Team member name (reuses the existing "Kevin Su" shared string elsewhere
# in the sheet, e.g. B75).
$ws.Cells.Item(76, 2).Value = "Kevin Su"

# DATE - 2025-01-19 (Excel serial date 45676)
$ws.Cells.Item(76, 3).Value = 45676

# TIME IN - 2:00 PM, TIME OUT - 4:30 PM (fractions of a 24h day)
$ws.Cells.Item(76, 4).Value = 0.58333333333333337
$ws.Cells.Item(76, 5).Value = 0.6875

# TOTAL - 2.5 hours
$ws.Cells.Item(76, 6).Value = 0.10416666666666667

# DESCRIPTION
$ws.Cells.Item(76, 7).Value = "Worked on mmWave Radar PCB schematic. Work included drawing components and connecting components."

# Match the row height Excel computes for this wrapped description text.
$ws.Rows.Item(76).RowHeight = 189

# Bring the newly-filled row into view and move the active selection to it,
# matching the saved view state (scrolled to row 76, selection on I76).
$ws.Activate()
$ws.Range("I76").Select() | Out-Null
